# "First Invalid login commit"
#
# Adds a new "InvalidLogin" worksheet (after the existing "ValidLogin" sheet)
# with a small data-driven test table of bad username/password combos and
# the expected failure message, then makes the new sheet the active tab.

$wb = $excel.ActiveWorkbook

# Existing sheet becomes the anchor so the new sheet lands right after it
# (Worksheets.Add defaults to inserting *before* the active sheet).
$validLogin = $wb.Worksheets.Item(1)

$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $validLogin)
$ws.Name = "InvalidLogin"

# Fill the data column-by-column (A then B then C) so that the shared-string
# table is populated in the same left-to-right, top-down order the strings
# were first authored in.
$ws.Range("A1").Value = "Username"
$ws.Range("A2").Value = "abcd"

$ws.Range("B1").Value = "Password"
$ws.Range("B2").Value = "xyz"

$ws.Range("C1").Value = "FailMsg"
$ws.Range("C2").Value = "Err Msg is Not Dispalyed"

$ws.Range("A3").Value = "admin"
$ws.Range("B3").Value = "damager"
$ws.Range("C3").Value = "Err Msg is Not Dispalyed"

$ws.Range("A4").Value = "admin"
$ws.Range("C4").Value = "Err Msg is Not Dispalyed"

$ws.Range("B5").Value = "manager"
$ws.Range("C5").Value = "Err Msg is Not Dispalyed"

# Best-fit the two text columns that hold the password / message values.
$ws.Columns.Item(2).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(3).EntireColumn.AutoFit() | Out-Null

# Select C5, zoom in, and make this new sheet the active tab - matching the
# saved view state of the edited workbook.
$ws.Range("C5").Select() | Out-Null
$excel.ActiveWindow.Zoom = 220
